# "coba no 22 ni :)"
# Java sheet: fix up row 17 (now just "coba" / "adeocaremiwowow"), and add
# new error-bank entries 18-22 as separate rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Java")
$ws.Activate()

# --- Row 20 (entry #17): shorten text, drop the trailing columns ---
$ws.Range("C20").Value = "coba"
$ws.Range("D20").Value = "adeocaremiwowow"
$ws.Range("E20").Clear()
$ws.Range("F20").Clear()
$ws.Range("G20").Clear()

# --- Row 21 (entry #18) ---
$ws.Range("B21").Value = 18
$ws.Range("B21").HorizontalAlignment = -4131
$ws.Range("C21").Value = "baru nih"
$ws.Range("C21").HorizontalAlignment = -4131
$ws.Range("F21").Value = "Rahmad"
$ws.Range("F21").HorizontalAlignment = -4131
$ws.Range("G21").Value = 43557
$ws.Range("G21").NumberFormat = "mm-dd-yy"

# --- Row 22 (entry #19) ---
$ws.Range("B22").Value = 19
$ws.Range("B22").HorizontalAlignment = -4131
$ws.Range("C22").Value = "Ini baru juga"
$ws.Range("C22").HorizontalAlignment = -4131

# --- Row 23 (entry #20) ---
$ws.Range("B23").Value = 20
$ws.Range("B23").HorizontalAlignment = -4131
$ws.Range("C23").Value = "ini ke 20"
$ws.Range("C23").HorizontalAlignment = -4131
$ws.Range("D23").Value = "khulqi"
$ws.Range("D23").HorizontalAlignment = -4131

# --- Row 24 (entry #21) ---
$ws.Range("B24").Value = 21
$ws.Range("B24").HorizontalAlignment = -4131
$ws.Range("C24").Value = "coba ke 21"
$ws.Range("C24").HorizontalAlignment = -4131
$ws.Range("D24").Value = "Oka"
$ws.Range("D24").HorizontalAlignment = -4131

# --- Row 25 (entry #22) ---
$ws.Range("B25").Value = 22
$ws.Range("B25").HorizontalAlignment = -4131
$ws.Range("C25").Value = "ke 22 ini"
$ws.Range("C25").HorizontalAlignment = -4131
$ws.Range("D25").Value = "fred"
$ws.Range("D25").HorizontalAlignment = -4131

# --- View state: scrolled down, D25 selected ---
$ws.Range("D25").Select()
$excel.ActiveWindow.ScrollRow = 11
$excel.ActiveWindow.ScrollColumn = 1
